$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text formatting so values
# like '0.9990' or '10.80' are not reinterpreted as numbers and
# lose trailing zeros.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '27.685.65'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '1.900.49'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.88%  '
$ws.Range('D5').Value = '311.77'
$ws.Range('E5').Value = '  -1.48%  '
$ws.Range('D6').Value = '0.9989'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('D7').Value = '0.5016'
$ws.Range('E7').Value = '  +3.89%  '
$ws.Range('D8').Value = '0.3767'
$ws.Range('E8').Value = '  -1.27%  '
$ws.Range('D9').Value = '0.07237'
$ws.Range('E9').Value = '  -1.67%  '
$ws.Range('D10').Value = '20.99'
$ws.Range('E10').Value = '  +1.39%  '
$ws.Range('D11').Value = '0.8904'
$ws.Range('E11').Value = '  -4.14%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.993.35'
$ws.Range('E12').Value = '  +4.73%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.07616'
$ws.Range('E13').Value = '  -1.62%  '
$ws.Range('D14').Value = '5.436'
$ws.Range('E14').Value = '  -0.83%  '
$ws.Range('D15').Value = '91.58'
$ws.Range('E15').Value = '  -0.09%  '
$ws.Range('D16').Value = '0.9990'
$ws.Range('D17').Value = '0.000008756'
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').Value = '27.736.72'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('D20').Value = '14.46'
$ws.Range('E20').Value = '  -1.15%  '
$ws.Range('D21').Value = '5.136'
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('D22').Value = '2.116.06'
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').Value = '10.80'
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').Value = '6.565'
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('D25').Value = '153.19'
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('D26').Value = '1.848'
$ws.Range('E26').Value = '  -3.71%  '
$ws.Range('D27').Value = '2.179'
$ws.Range('E27').Value = '  +2.40%  '
$ws.Range('D28').Value = '18.26'
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('D29').Value = '114.67'
$ws.Range('E29').Value = '  -2.12%  '
$ws.Range('D30').Value = '4.824'
$ws.Range('E30').Value = '  -2.48%  '
$ws.Range('D31').Value = '0.08924'
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('D32').Value = '3.180'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '4.788'
$ws.Range('E33').Value = '  +3.04%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = '1.228'
$ws.Range('E34').Value = '  -1.71%  '
$ws.Range('D35').Value = '0.7825'
$ws.Range('E35').Value = '  +2.20%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '2.623'
$ws.Range('E36').Value = '  +3.86%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.02083'
$ws.Range('E37').Value = '  +1.95%  '
$ws.Range('D38').Value = '3.053'
$ws.Range('E38').Value = '  +1.85%  '
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').Value = '0.5487'
$ws.Range('E40').Value = '  +0.29%  '
$ws.Range('D41').Value = '0.05276'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = '6.715'
$ws.Range('E42').Value = '  -3.43%  '
$ws.Range('D43').Value = '113.46'
$ws.Range('E43').Value = '  +3.86%  '
$ws.Range('D44').Value = '8.437'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('D45').Value = '0.1510'
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('D46').Value = '0.4767'
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('D47').Value = '10.45'
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('D48').Value = '0.9987'
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').Value = '1.611'
$ws.Range('E49').Value = '  -2.30%  '
$ws.Range('D50').Value = '66.53'
$ws.Range('E50').Value = '  -1.60%  '
$ws.Range('D51').Value = '0.05995'
$ws.Range('E51').Value = '  -1.49%  '
